$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://cdn.jsdelivr.net/npm/emoji-datasource-apple/img/apple/64/1f4c8.png"
$ws.Range("B2").Value = 10000
$ws.Range("C2").Value = "Trading Loss"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2025-04-24"
